$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date value for A78 (2025-05-01 -> serial 45778), formatted as a date like A77
$ws.Range("A77").Copy()
$ws.Range("A78").PasteSpecial(-4122)
$ws.Range("A78").Value = 45778

$ws.Range("B78").Value = -0.541
$ws.Range("C78").Value = -0.186
$ws.Range("D78").Value = 0.216
$ws.Range("E78").Value = 0.248
$ws.Range("F78").Value = 0.076
$ws.Range("G78").Value = 79.71
